# Babine fence counts update for 10/7 (rows 94-97, dates 45566-45569)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 94 (2024-10-07 -> serial 45566) ----
$ws.Range("B94").Value = 60
$ws.Range("C94").Value = 15
$ws.Range("D94").Value = 17
$ws.Range("E94").Value = 59
$ws.Range("F94").Value = 4
$ws.Range("G94").Value = 5
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Formula = "=M93+B94"
$ws.Range("N94").Formula = "=J94+L94"
$ws.Range("O94").Formula = "=B94+I94+N94"
$ws.Range("P94").Formula = "=P93+O94"
$ws.Range("Q94").Formula = "=C94+K94"
$ws.Range("R94").Formula = "=Q94+R93"
$ws.Range("S94").Value = 897

# ---- Row 95 (45567) ----
$ws.Range("B95").Value = 62
$ws.Range("C95").Value = 10
$ws.Range("D95").Value = 8
$ws.Range("E95").Value = 23
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").Formula = "=M94+B95"
$ws.Range("N95").Formula = "=J95+L95"
$ws.Range("O95").Formula = "=B95+I95+N95"
$ws.Range("P95").Formula = "=P94+O95"
$ws.Range("Q95").Formula = "=C95+K95"
$ws.Range("R95").Formula = "=Q95+R94"
$ws.Range("S95").Value = 898

# ---- Row 96 (45568) ----
$ws.Range("B96").Value = 60
$ws.Range("C96").Value = 5
$ws.Range("D96").Value = 15
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Formula = "=M95+B96"
$ws.Range("N96").Formula = "=J96+L96"
$ws.Range("O96").Formula = "=B96+I96+N96"
$ws.Range("P96").Formula = "=P95+O96"
$ws.Range("Q96").Formula = "=C96+K96"
$ws.Range("R96").Formula = "=Q96+R95"
$ws.Range("S96").Value = 899

# ---- Row 97 (45569) ----
$ws.Range("B97").Value = 48
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 15
$ws.Range("E97").Value = 11
$ws.Range("F97").Value = 1
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Formula = "=M96+B97"
$ws.Range("N97").Formula = "=J97+L97"
$ws.Range("O97").Formula = "=B97+I97+N97"
$ws.Range("P97").Formula = "=P96+O97"
$ws.Range("Q97").Formula = "=C97+K97"
$ws.Range("R97").Formula = "=Q97+R96"
$ws.Range("S97").Value = 900

# Recalculate so all dependent cells (U/V columns cascading through row 108, etc.) refresh.
$excel.Calculate()

# ---- View-state tweaks (best effort) ----
$ws.Range("Q100").Select()
